$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51 (columns B: Coin, C: Link, D: Price, E: Volume(1h))
$data = @(
    @{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "28.715.94"; E = "  +1.83%  " },
    @{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "1.892.57"; E = "  +0.65%  " },
    @{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.032"; E = "  +2.46%  " },
    @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "320.04"; E = "  +1.75%  " },
    @{ Row = 6; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.028"; E = "  +2.01%  " },
    @{ Row = 7; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.5181"; E = "  +0.90%  " },
    @{ Row = 8; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.4007"; E = "  +2.66%  " },
    @{ Row = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.08415"; E = "  +0.41%  " },
    @{ Row = 10; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "1.127"; E = "  +0.46%  " },
    @{ Row = 11; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "42.45"; E = "  +2.13%  " },
    @{ Row = 12; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "6.320"; E = "  +1.47%  " },
    @{ Row = 13; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "20.70"; E = "  +0.10%  " },
    @{ Row = 14; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "7.303"; E = "  +0.71%  " },
    @{ Row = 15; B = "BinanceUSD"; C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D = "1.032"; E = "  +2.35%  " },
    @{ Row = 16; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "1.789.14"; E = "  -6.16%  " },
    @{ Row = 17; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.00001119"; E = "  +1.51%  " },
    @{ Row = 18; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "91.98"; E = "  +0.98%  " },
    @{ Row = 19; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.06831"; E = "  +1.84%  " },
    @{ Row = 20; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "17.90"; E = "  +0.53%  " },
    @{ Row = 21; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.028"; E = "  +2.05%  " },
    @{ Row = 22; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "6.039"; E = "  -0.03%  " },
    @{ Row = 23; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "28.755.28"; E = "  +1.79%  " },
    @{ Row = 24; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "11.25"; E = "  +1.00%  " },
    @{ Row = 25; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "2.294"; E = "  +1.06%  " },
    @{ Row = 26; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "163.12"; E = "  +2.39%  " },
    @{ Row = 27; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "21.02"; E = "  +1.96%  " },
    @{ Row = 28; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "2.023.35"; E = "  -3.36%  " },
    @{ Row = 29; B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "2.410"; E = "  -2.81%  " },
    @{ Row = 30; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "128.91"; E = "  +2.95%  " },
    @{ Row = 31; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.1064"; E = "  +0.29%  " },
    @{ Row = 32; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.054"; E = "  +1.55%  " },
    @{ Row = 33; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "5.890"; E = "  +0.46%  " },
    @{ Row = 34; B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "3.672"; E = "  +1.73%  " },
    @{ Row = 35; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.02459"; E = "  -0.17%  " },
    @{ Row = 36; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.06575"; E = "  -0.07%  " },
    @{ Row = 37; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "9.251"; E = "  -3.81%  " },
    @{ Row = 38; B = "Algorand"; C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D = "0.2204"; E = "  +0.63%  " },
    @{ Row = 39; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D = "1.284"; E = "  +4.84%  " },
    @{ Row = 40; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "1.201"; E = "  +0.14%  " },
    @{ Row = 41; B = "TheSandbox"; C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D = "0.6503"; E = "  -0.18%  " },
    @{ Row = 42; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "5.063"; E = "  +1.21%  " },
    @{ Row = 43; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "11.29"; E = "  -0.27%  " },
    @{ Row = 44; B = "Decentraland"; C = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D = "0.6094"; E = "  -1.00%  " },
    @{ Row = 45; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "13.12"; E = "  +0.28%  " },
    @{ Row = 46; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "3.761"; E = "  +2.13%  " },
    @{ Row = 47; B = "WEMIXTOKEN"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "1.237"; E = "  -3.93%  " },
    @{ Row = 48; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "2.020"; E = "  +0.29%  " },
    @{ Row = 49; B = "EOS"; C = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D = "1.221"; E = "  -0.57%  " },
    @{ Row = 50; B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "122.64"; E = "  +1.14%  " },
    @{ Row = 51; B = "Cronos"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "0.06906"; E = "  -0.13%  " }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    # Force column D (Price) to be stored as text so that values like
    # "28.715.94", "6.320", "320.04" keep their exact textual form
    # instead of being auto-converted into numbers.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
